# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 59 (pushing the existing rows 59-153
# down to 60-154) and populate it with the new week's observation for
# Poroto verde at Macroferia Regional de Talca.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 59:153 down to 60:154, carrying formatting from row 59 (date
# style on column D, etc.) to the newly inserted row, same as Excel's
# default "Insert" behavior (shift cells down, inherit formats from above).
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with this week's record.
$ws.Range("A59").Value = 5
$ws.Range("B59").Value = "Macroferia Regional de Talca"
$ws.Range("C59").Value = "Maule"
$ws.Range("D59").Value = 44638
$ws.Range("E59").Value = 7
$ws.Range("F59").Value = 100112031
$ws.Range("G59").Value = "Poroto verde"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 100
$ws.Range("K59").Value = 25000
$ws.Range("L59").Value = 25000
$ws.Range("M59").Value = 25000
$ws.Range("N59").Value = "$/saco 25 kilos"
$ws.Range("O59").Value = "Región del Maule"
$ws.Range("P59").Value = 1000
$ws.Range("Q59").Value = 25
$ws.Range("R59").Value = "Hortaliza"

# Keep the inserted date cell formatted the same way as the rest of
# column D (yyyy-mm-dd hh:mm:ss).
$ws.Range("D59").NumberFormat = "YYYY-MM-DD HH:MM:SS"
